# Apply the "yearly.xlsx" cash-flow update described by the commit.
# The workbook has a single sheet ("Overview"). All edits are plain
# literal values (no formulas involved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Header date labels (row 9), columns G/H ---------------------------
$ws.Range("G9").Value = "1402-04-14 (9)"
$ws.Range("H9").Value = "1402-04-14 (2)"

# --- Operating activities section (rows 12-17) --------------------------
$ws.Range("G12").Value = 45835541
$ws.Range("H12").Value = 42544775

$ws.Range("G14").Value = 45820121
$ws.Range("H14").Value = 42544775

$ws.Range("G16").Value = 27181
$ws.Range("H16").Value = 2091238

$ws.Range("H17").Value = -25934829

# --- Investing activities section (rows 29, 32) -------------------------
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0

$ws.Range("G32").Value = -1267910
$ws.Range("H32").Value = -22749910

# --- Financing activities section (row 36) -------------------------------
# F36 switches from the textual placeholder "-" to a literal numeric 0.
$ws.Range("F36").Value = 0
